$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 holds the "Enterprises density (per 1000 people)" figures for
# Micro (B13), SMEs (C13) and MSMEs (D13). They are stored as text
# (shared-string) values "6.9", "1.9", "8.8" and need to become the more
# precise "6.89", "1.94", "8.83" while staying plain text cells (not
# becoming numeric cells).
#
# Temporarily force a text number format before assigning the value so
# Excel does not auto-convert the numeric-looking string into a number,
# then restore the original General format/Normal style so the cells'
# formatting matches what they had before the edit.

$target = $ws.Range("B13:D13")
$target.NumberFormat = "@"

$ws.Range("B13").Value = "6.89"
$ws.Range("C13").Value = "1.94"
$ws.Range("D13").Value = "8.83"

$target.NumberFormat = "General"
$target.Style = "Normal"
